# Generate Report for Handback
# a.md has been handed back and is now in sync with en-US; update the
# status/handback metadata on all three sheets to reflect that, while
# b.md remains "In Translation" (its cells are untouched in content).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: a.md (row 2) is now handed back / in sync
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Columns.Item(5).ColumnWidth = 29.14
$overview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------
# zh-cn sheet: a.md (row 2) handback info
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("L2").Value = "2017-02-09 08:19:59"
$zhcn.Range("M2").Value = "TestHandback_201702090419"
$zhcn.Range("R2").Value = ""
$zhcn.Columns.Item(13).ColumnWidth = 27.15

# ---------------------------------------------------------------
# de-de sheet: a.md (row 2) handback info
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("L2").Value = "2017-02-09 08:20:26"
$dede.Range("M2").Value = "TestHandback_201702090419"
$dede.Range("R2").Value = ""
$dede.Columns.Item(13).ColumnWidth = 27.15
